# análises - revisão da exportação
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header style (bold + border, currently on A1:D1) across the
# new header columns E1:K1 before touching any values.
$ws.Range("A1:D1").Copy()
$ws.Range("E1:K1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the old data rows (rows 2-10) entirely - the new table only needs
# rows 2-3.
$ws.Rows("4:10").Delete()

# New header row
$headers = @("origem","total","aon","aon_sucesso","aon_falha","flex","flex_sucesso","flex_falha","sub","sub_sucesso","sub_falha")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# New data rows
$data = @(
    @("apoia.se", 632, 0, 0, 0, 5, 0, 5, 627, 137, 490),
    @("catarse", 2855, 1335, 830, 505, 1463, 1383, 80, 57, 15, 42)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
